$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 32   Number  41"
$ws.Range("C9").Value = "Report Covering the Week  10/6/2025  Through  10/12/2025"

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 2
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("H15").Value = 100
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = 45.454545454545
$ws.Range("L15").Value = 77.777777777777
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = -44.827586206896

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 80
$ws.Range("I16").Value = 78
$ws.Range("J16").Value = 110
$ws.Range("K16").Value = -29.090909090909
$ws.Range("L16").Value = -22.772277227722
$ws.Range("M16").Value = -48.684210526315
$ws.Range("N16").Value = -88.936170212766

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 6.25
$ws.Range("I17").Value = 144
$ws.Range("J17").Value = 122
$ws.Range("K17").Value = 18.032786885245
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 77.777777777777
$ws.Range("N17").Value = -63.358778625954

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -21.428571428571
$ws.Range("I18").Value = 141
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = 42.424242424242
$ws.Range("L18").Value = -0.704225352112
$ws.Range("M18").Value = 35.576923076923
$ws.Range("N18").Value = -83.293838862559

# Row 19
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = -14.545454545454
$ws.Range("I19").Value = 425
$ws.Range("J19").Value = 407
$ws.Range("K19").Value = 4.422604422604
$ws.Range("L19").Value = 5.985037406483
$ws.Range("M19").Value = 2.657004830917
$ws.Range("N19").Value = -50.234192037470

# Row 20
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 32
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = -43.859649122807
$ws.Range("L20").Value = -60.493827160493
$ws.Range("M20").Value = 3.225806451612
$ws.Range("N20").Value = -95.928753180661

# Row 21
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -21.212121212121
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -10.679611650485
$ws.Range("I21").Value = 837
$ws.Range("J21").Value = 806
$ws.Range("K21").Value = 3.846153846153
$ws.Range("L21").Value = -4.994324631101
$ws.Range("M21").Value = 5.415617128463
$ws.Range("N21").Value = -76.948499036078

# Row 22
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = -18.518518518518

# Row 23
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 18.181818181818
$ws.Range("I23").Value = 93
$ws.Range("J23").Value = 81
$ws.Range("K23").Value = 14.814814814814
$ws.Range("L23").Value = 6.896551724137
$ws.Range("M23").Value = 50

# Row 24
$ws.Range("C24").Value = 47
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 67.857142857142
$ws.Range("F24").Value = 175
$ws.Range("G24").Value = 126
$ws.Range("H24").Value = 38.888888888888
$ws.Range("I24").Value = 1374
$ws.Range("J24").Value = 937
$ws.Range("K24").Value = 46.638207043756
$ws.Range("L24").Value = 11.075181891673
$ws.Range("M24").Value = 62.219598583234

# Row 25
$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 70
$ws.Range("F25").Value = 123
$ws.Range("G25").Value = 85
$ws.Range("H25").Value = 44.705882352941
$ws.Range("I25").Value = 919
$ws.Range("J25").Value = 524
$ws.Range("K25").Value = 75.381679389313
$ws.Range("L25").Value = 7.737397420867

# Row 26
$ws.Range("C26").Value = 5
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = 47.619047619047
$ws.Range("I26").Value = 248
$ws.Range("J26").Value = 246
$ws.Range("K26").Value = 0.813008130081
$ws.Range("L26").Value = 7.826086956521
$ws.Range("M26").Value = -6.415094339622

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = 0

# Row 29
$ws.Range("L29").Value = -57.142857142857

# Row 30
$ws.Range("L30").Value = -66.666666666666

# Row 31
$ws.Range("D31").Value = 2
$ws.Range("D31").NumberFormat = '#,##0'
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 20
$ws.Range("K31").Value = -55

# --- Cells changing from numeric to shared-string text ("0" / "***.*") ---
# Template source cells (untouched elsewhere): C14 = "0" (style 13), E14 = "***.*" (style 13)
$ws.Range("C20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("C31").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C31").PasteSpecial(-4122)

$excel.CutCopyMode = 0
